# Regenerate the localization-status report: items that have moved out of
# the handoff queue and are now being translated get their status updated
# from "Ready for handoff" to "In Translation" on every sheet that surfaces
# status (the Overview rollup columns per-locale, and each locale's own
# Status column). Excel re-fits the (now shorter) status column after the
# text changes, so we re-apply AutoFit to those columns as well.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: per-locale status rollup columns (zh-cn = E, de-de = F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = $newStatus
$wsOverview.Columns("E:F").AutoFit()
$wsOverview.Columns("E:F").ColumnWidth = 12.5

# --- zh-cn sheet: Status column (C) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2:C3").Value = $newStatus
$wsZh.Columns("C:C").AutoFit()
$wsZh.Columns("C:C").ColumnWidth = 12.5

# --- de-de sheet: Status column (C) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2:C3").Value = $newStatus
$wsDe.Columns("C:C").AutoFit()
$wsDe.Columns("C:C").ColumnWidth = 12.5
